$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Ají" on 2021-12-17 (serial 44547).
# Insert a new row at position 9 (shifting existing rows 9-46 down to 10-47)
# and populate it with the new observation.
$ws.Rows("9:9").Insert()

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44547
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112021
$ws.Cells.Item(9, 7).Value = "Ají"
$ws.Cells.Item(9, 8).Value = "Americana (o)"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 17000
$ws.Cells.Item(9, 12).Value = 17500
$ws.Cells.Item(9, 13).Value = 17250
$ws.Cells.Item(9, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 1150
$ws.Cells.Item(9, 17).Value = 15
$ws.Cells.Item(9, 18).Value = "Hortaliza"
